$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (row 2 through row 15) holds a "Förändrad" (changed/updated) date
# that was bumped forward by one day, from serial 45189 (2023-09-20) to
# serial 45190 (2023-09-21). Update each cell in that range.
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45190
}
